# chore: publish terminology IG 2.0.0 (#33)
#
# 1. Metadata!Version 1.8.1 -> 1.0.1, Metadata!Date 2022-09-01 -> 2025-09-22
# 2. Rename existing "Concepts" sheet to "Properties" and replace its
#    contents with the CodeSystem properties table (Code/Uri/Description/Type).
# 3. Add a fresh "Concepts" sheet (after "Properties") holding the table that
#    used to live in the old "Concepts" sheet (Level/Code/Display/Definition).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- helper: force a literal TEXT value into a cell, bypassing Excel's
# automatic number/date detection (which would otherwise turn
# "2025-09-22" into a date serial, or "1" into a quote-prefixed number
# and, either way, mint a brand-new cell style). We build the literal via
# a throwaway text-formula cell, copy its computed value (value-only
# paste keeps the destination's existing style untouched), then wipe the
# scratch cell again.
function Set-LiteralText($Worksheet, $Address, $Text) {
    $scratch = $Worksheet.Range("ZZ1")
    $escaped = $Text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $Worksheet.Range($Address).PasteSpecial(4163)
    $scratch.Value = $null
}

# 1. Metadata updates -------------------------------------------------
$ws1.Range("B3").Value = "1.0.1"
Set-LiteralText $ws1 "B8" "2025-09-22"

# 2. Turn the old "Concepts" sheet into "Properties" -------------------
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Name = "Properties"
$props = $wb.Worksheets.Item("Properties")

# Extend formatting (header style row1, data style row2) down to row3
# before touching values, so the new row inherits the right cell style.
$props.Range("A2:D2").Copy()
$props.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$props.Range("A1").Value = "Code"
$props.Range("B1").Value = "Uri"
$props.Range("C1").Value = "Description"
$props.Range("D1").Value = "Type"

$props.Range("A2").Value = "status"
$props.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$props.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$props.Range("D2").Value = "code"

$props.Range("A3").Value = "effectiveDate"
$props.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$props.Range("C3").Value = "The date at which the concept status was last changed."
$props.Range("D3").Value = "dateTime"

# 3. Add the new "Concepts" sheet (right after "Properties") -----------
$newSheet = $wb.Worksheets.Add($null, $props)
$newSheet.Name = "Concepts"
$newConcepts = $wb.Worksheets.Item("Concepts")

# Borrow the header/data styles from the "Properties" sheet so the new
# sheet's cells line up with s=1 (header) / s=2 (data) instead of
# getting no explicit style at all.
$props.Range("A1:D1").Copy()
$newConcepts.Range("A1:D1").PasteSpecial(-4122)
$props.Range("A2:D2").Copy()
$newConcepts.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newConcepts.Range("A1").Value = "Level"
$newConcepts.Range("B1").Value = "Code"
$newConcepts.Range("C1").Value = "Display"
$newConcepts.Range("D1").Value = "Definition"

Set-LiteralText $newConcepts "A2" "1"
$newConcepts.Range("B2").Value = "other"
$newConcepts.Range("C2").Value = "Other encounter class"

# Leave the active tab as it was originally ("Metadata" / activeTab=0).
$ws1.Activate()
